# Add a new "T4: 18/3/2020" column (G) to the COVID19 history sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("G1").Value = "T4: 18/3/2020"

# New column data (rows 2-19), one value per department row
$gValues = @(1, 2, 0, 0, 0, 2, 0, 7, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}

# Totals row: sum of the new column
$ws.Range("G20").Formula = "=SUM(G2:G19)"

# Carry over the "El Paraiso" row's highlighted/underlined format from F8 to G8
$ws.Range("F8").Copy()
$ws.Range("G8").PasteSpecial(-4122)

# Move the active selection, matching where the editor's cursor ended up
$ws.Range("G10").Select()
